$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# format first, so Excel stores them as the exact original text (preserving
# things like leading/trailing zeros, e.g. '0.140') instead of re-parsing them
# into floating point numbers. (Looping individually since a Union range's
# NumberFormat setter only reliably applies to the first area.)
$textCells = @('D4', 'D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D17', 'D19', 'D21', 'D23', 'D24', 'D26', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D38', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D49', 'D50', 'D51')
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '51.935.49'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '2.834.96'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '356.74'
$ws.Range("E5").Value = '  +4.08%  '
$ws.Range("D6").Value = '112.59'
$ws.Range("E6").Value = '  -2.66%  '
$ws.Range("D7").Value = '0.567'
$ws.Range("E7").Value = '  +3.40%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.603'
$ws.Range("E9").Value = '  +4.25%  '
$ws.Range("D10").Value = '41.16'
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '20.07'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").Value = '3.276.14'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '2.839.62'
$ws.Range("E16").Value = '  +2.12%  '
$ws.Range("D17").Value = '0.932'
$ws.Range("E17").Value = '  +6.37%  '
$ws.Range("D18").Value = '51.846.16'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").Value = '7.54'
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("D21").Value = '13.45'
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").Value = '0.0₃0993'
$ws.Range("E22").Value = '  +1.43%  '
$ws.Range("D23").Value = '70.08'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '269.46'
$ws.Range("E24").Value = '  -2.75%  '
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("D26").Value = '27.07'
$ws.Range("E26").Value = '  +1.45%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("B31").Value = 'VeChain'
$ws.Range("C31").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D31").Value = '0.0476'
$ws.Range("E31").Value = '  +24.37%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '35.50'
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '52.61'
$ws.Range("E33").Value = '  +4.93%  '
$ws.Range("D34").Value = '5.91'
$ws.Range("E34").Value = '  +3.15%  '
$ws.Range("D35").Value = '5.43'
$ws.Range("E35").Value = '  +9.37%  '
$ws.Range("D36").Value = '0.0846'
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '3.29'
$ws.Range("E38").Value = '  +2.10%  '
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("D41").Value = '0.117'
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '2.55'
$ws.Range("E42").Value = '  -4.74%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '23.25'
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '124.04'
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").Value = '2.28'
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '3.38'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.100.86'
$ws.Range("E47").Value = '  +1.39%  '
$ws.Range("D49").Value = '5.99'
$ws.Range("E49").Value = '  +7.74%  '
$ws.Range("D50").Value = '0.977'
$ws.Range("E50").Value = '  +9.06%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '9.08'
$ws.Range("E51").Value = '  +2.62%  '
